# Apply the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.752.07"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.400.17"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "'" + "563.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").Value = "'" + "141.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").Value = "2.406.99"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "'" + "0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "'" + "5.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "'" + "0.342"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "'" + "26.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "2.836.22"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "60.686.16"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "2.405.86"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'" + "8.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.15%  "
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "'" + "324.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D23").Value = "'" + "6.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'" + "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("D26").Value = "'" + "65.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'" + "564.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("D28").Value = "'" + "8.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.37%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "0.0₃0936"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").Value = "'" + "8.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").Value = "'" + "0.131"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("D35").Value = "'" + "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  +3.53%  "
$ws.Range("D37").Value = "'" + "152.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "'" + "0.371"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "'" + "4.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").Value = "'" + "18.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'" + "5.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'" + "2.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.29%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'" + "41.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'" + "1.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").Value = "'" + "141.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").Value = "'" + "3.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "'" + "0.591"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "'" + "0.0508"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'" + "19.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.34%  "
